$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GHGs")

$data = @(
  @("PWRCCGCRIXX00", "CO2CRI", 0.1283, "MtCO2eq/PJ"),
  @("PWRCCGCRIXX01", "CO2CRI", 0.1006, "MtCO2eq/PJ"),
  @("PWRCCGPANXX00", "CO2PAN", 0.1232, "MtCO2eq/PJ"),
  @("PWRCCGPANXX01", "CO2PAN", 0.1006, "MtCO2eq/PJ"),
  @("PWRCCSCRIXX01", "CO2CRI", 0.0273, "MtCO2eq/PJ"),
  @("PWRCCSPANXX01", "CO2PAN", 0.0273, "MtCO2eq/PJ"),
  @("PWRCOACRIXX01", "CO2CRI", 0.2724, "MtCO2eq/PJ"),
  @("PWRCOAPANXX01", "CO2PAN", 0.2771, "MtCO2eq/PJ"),
  @("PWRCOGCRIXX01", "CO2CRI", 0.2724, "MtCO2eq/PJ"),
  @("PWRCOGPANXX01", "CO2PAN", 0.2724, "MtCO2eq/PJ"),
  @("PWROCGCRIXX00", "CO2CRI", 0.1504, "MtCO2eq/PJ"),
  @("PWROCGCRIXX01", "CO2CRI", 0.1437, "MtCO2eq/PJ"),
  @("PWROCGPANXX00", "CO2PAN", 0.1504, "MtCO2eq/PJ"),
  @("PWROCGPANXX01", "CO2PAN", 0.1437, "MtCO2eq/PJ"),
  @("PWROILCRIXX01", "CO2CRI", 0.2021, "MtCO2eq/PJ"),
  @("PWROILPANXX01", "CO2PAN", 0.2077, "MtCO2eq/PJ"),
  @("PWROTHCRIXX01", "CO2CRI", 0.0503, "MtCO2eq/PJ"),
  @("PWROTHPANXX01", "CO2PAN", 0.0503, "MtCO2eq/PJ"),
  @("PWRPETCRIXX01", "CO2CRI", 0.217, "MtCO2eq/PJ"),
  @("PWRPETPANXX01", "CO2PAN", 0.217, "MtCO2eq/PJ"),
  @("PWRWASCRIXX01", "CO2CRI", 0.0879, "MtCO2eq/PJ"),
  @("PWRWASPANXX01", "CO2PAN", 0.0879, "MtCO2eq/PJ"),
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = 2 + $i
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}
